$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the English title of the "Dragon Continent" row (C7):
# "Dragon Contient: Hyper T&T World Guide" -> "Dragon Continent: Hyper T&T World Guide"
$ws.Range("C7").Value = "Dragon Continent: Hyper T&T World Guide"

# Update the active selection to B8 (as recorded in the saved workbook state)
$ws.Range("B8").Select()
